$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1942.5238
$ws.Cells.Item(18, 9).Value = 2027.5
$ws.Cells.Item(18, 10).Value = 1432.6666
$ws.Cells.Item(18, 11).Value = 2027.5
$ws.Cells.Item(18, 12).Value = 1432.6666
$ws.Cells.Item(18, 13).Value = -1743.5
$ws.Cells.Item(18, 14).Value = -2000.6666

$ws.Cells.Item(86, 8).Value = 9023.177
$ws.Cells.Item(86, 9).Value = 8489.5
$ws.Cells.Item(86, 11).Value = 8489.5
$ws.Cells.Item(86, 13).Value = -7366.5

$ws.Cells.Item(89, 8).Value = 9023.177
$ws.Cells.Item(89, 9).Value = 8489.5
$ws.Cells.Item(89, 11).Value = 42447.5
$ws.Cells.Item(89, 13).Value = -36831.5

$ws.Cells.Item(98, 8).Value = 1506.5652
$ws.Cells.Item(98, 9).Value = 1435.8096
$ws.Cells.Item(98, 10).Value = 2249.5
$ws.Cells.Item(98, 11).Value = 1435.8096
$ws.Cells.Item(98, 12).Value = 2249.5
$ws.Cells.Item(98, 13).Value = 62.19039999999995
$ws.Cells.Item(98, 14).Value = -5245.5

$ws.Cells.Item(101, 8).Value = 100000720
$ws.Cells.Item(101, 10).Value = 1000
$ws.Cells.Item(101, 12).Value = 3000
$ws.Cells.Item(101, 14).Value = -6244

$ws.Cells.Item(106, 8).Value = 41676840
$ws.Cells.Item(106, 9).Value = 71443300
$ws.Cells.Item(106, 11).Value = 71443300
$ws.Cells.Item(106, 13).Value = -71442669

$ws.Cells.Item(112, 8).Value = 2847.0908

$ws.Cells.Item(115, 8).Value = 2083.1667
$ws.Cells.Item(115, 10).Value = 3000
$ws.Cells.Item(115, 12).Value = 9000
$ws.Cells.Item(115, 14).Value = -12134

$ws.Cells.Item(116, 8).Value = 6874.9287
$ws.Cells.Item(116, 9).Value = 3303.4
$ws.Cells.Item(116, 10).Value = 8859.111000000001
$ws.Cells.Item(116, 11).Value = 3303.4
$ws.Cells.Item(116, 12).Value = 8859.111000000001
$ws.Cells.Item(116, 13).Value = 138.5999999999999
$ws.Cells.Item(116, 14).Value = -15743.111

$ws.Cells.Item(121, 8).Value = 2817.2
$ws.Cells.Item(121, 10).Value = 2817.2
$ws.Cells.Item(121, 12).Value = 8451.599999999999
$ws.Cells.Item(121, 14).Value = -11945.6

$ws.Cells.Item(122, 8).Value = 1506.5652
$ws.Cells.Item(122, 9).Value = 1435.8096
$ws.Cells.Item(122, 10).Value = 2249.5
$ws.Cells.Item(122, 11).Value = 4307.4288
$ws.Cells.Item(122, 12).Value = 6748.5
$ws.Cells.Item(122, 13).Value = -1857.4288
$ws.Cells.Item(122, 14).Value = -11648.5

$ws.Cells.Item(127, 8).Value = 974.1111
$ws.Cells.Item(127, 9).Value = 783.375
$ws.Cells.Item(127, 10).Value = 2500
$ws.Cells.Item(127, 11).Value = 2350.125
$ws.Cells.Item(127, 12).Value = 7500
$ws.Cells.Item(127, 13).Value = 2609.875
$ws.Cells.Item(127, 14).Value = -17420

$ws.Cells.Item(132, 8).Value = 3107.5117
$ws.Cells.Item(132, 9).Value = 3414.2222
$ws.Cells.Item(132, 10).Value = 1530.1428
$ws.Cells.Item(132, 11).Value = 10242.6666
$ws.Cells.Item(132, 12).Value = 4590.428400000001
$ws.Cells.Item(132, 13).Value = -7712.6666
$ws.Cells.Item(132, 14).Value = -9650.428400000001

$ws.Cells.Item(137, 8).Value = 68332.92999999999
$ws.Cells.Item(137, 10).Value = 3032.3
$ws.Cells.Item(137, 12).Value = 9096.900000000001
$ws.Cells.Item(137, 14).Value = -14196.9

$ws.Cells.Item(138, 8).Value = 3655.4546
$ws.Cells.Item(138, 10).Value = 3679.5454
$ws.Cells.Item(138, 12).Value = 11038.6362
$ws.Cells.Item(138, 14).Value = -21318.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4902.95
$ws.Cells.Item(32, 9).Value = 4065.823
$ws.Cells.Item(32, 10).Value = 24994
$ws.Cells.Item(32, 11).Value = 4065.823
$ws.Cells.Item(32, 12).Value = 24994
$ws.Cells.Item(32, 13).Value = -3778.823
$ws.Cells.Item(32, 14).Value = -25568

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 78394.5
$ws.Cells.Item(140, 10).Value = 78131.60000000001
$ws.Cells.Item(140, 12).Value = 78131.60000000001
$ws.Cells.Item(140, 14).Value = -88491.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 7953
$ws.Cells.Item(58, 9).Value = 10184.75
$ws.Cells.Item(58, 10).Value = 4127.143
$ws.Cells.Item(58, 11).Value = 10184.75
$ws.Cells.Item(58, 12).Value = 4127.143
$ws.Cells.Item(58, 13).Value = -9981.75
$ws.Cells.Item(58, 14).Value = -4533.143

$ws.Cells.Item(86, 8).Value = 7872.0835
$ws.Cells.Item(86, 10).Value = 14317.167
$ws.Cells.Item(86, 12).Value = 14317.167
$ws.Cells.Item(86, 14).Value = -16563.167

$ws.Cells.Item(89, 8).Value = 7872.0835
$ws.Cells.Item(89, 10).Value = 14317.167
$ws.Cells.Item(89, 12).Value = 71585.83499999999
$ws.Cells.Item(89, 14).Value = -82817.83499999999

$ws.Cells.Item(99, 8).Value = 3676.25
$ws.Cells.Item(99, 9).Value = 3380.25
$ws.Cells.Item(99, 10).Value = 4120.25
$ws.Cells.Item(99, 11).Value = 3380.25
$ws.Cells.Item(99, 12).Value = 4120.25
$ws.Cells.Item(99, 13).Value = -1882.25
$ws.Cells.Item(99, 14).Value = -7116.25

$ws.Cells.Item(107, 8).Value = 2728.577
$ws.Cells.Item(107, 9).Value = 2379.318
$ws.Cells.Item(107, 10).Value = 4649.5
$ws.Cells.Item(107, 11).Value = 2379.318
$ws.Cells.Item(107, 12).Value = 4649.5
$ws.Cells.Item(107, 13).Value = -459.3180000000002
$ws.Cells.Item(107, 14).Value = -8489.5

$ws.Cells.Item(126, 8).Value = 3676.25
$ws.Cells.Item(126, 9).Value = 3380.25
$ws.Cells.Item(126, 10).Value = 4120.25
$ws.Cells.Item(126, 11).Value = 10140.75
$ws.Cells.Item(126, 12).Value = 12360.75
$ws.Cells.Item(126, 13).Value = -7670.75
$ws.Cells.Item(126, 14).Value = -17300.75

$ws.Cells.Item(136, 8).Value = 7953
$ws.Cells.Item(136, 9).Value = 10184.75
$ws.Cells.Item(136, 10).Value = 4127.143
$ws.Cells.Item(136, 11).Value = 30554.25
$ws.Cells.Item(136, 12).Value = 12381.429
$ws.Cells.Item(136, 13).Value = -28004.25
$ws.Cells.Item(136, 14).Value = -17481.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 1970.5714
$ws.Cells.Item(80, 9).Value = 1490
$ws.Cells.Item(80, 10).Value = 3172
$ws.Cells.Item(80, 11).Value = 4470
$ws.Cells.Item(80, 12).Value = 9516
$ws.Cells.Item(80, 13).Value = -3534
$ws.Cells.Item(80, 14).Value = -11388

$ws.Cells.Item(83, 8).Value = 1970.5714
$ws.Cells.Item(83, 9).Value = 1490
$ws.Cells.Item(83, 10).Value = 3172
$ws.Cells.Item(83, 11).Value = 13410
$ws.Cells.Item(83, 12).Value = 28548
$ws.Cells.Item(83, 13).Value = -8730
$ws.Cells.Item(83, 14).Value = -37908

$ws.Cells.Item(133, 8).Value = 3888.3333
$ws.Cells.Item(133, 9).Value = 3888.3333
$ws.Cells.Item(133, 11).Value = 11664.9999
$ws.Cells.Item(133, 13).Value = -6604.999899999999

$ws.Cells.Item(137, 8).Value = 5748.375
$ws.Cells.Item(137, 9).Value = 3749.6667
$ws.Cells.Item(137, 10).Value = 6414.6113
$ws.Cells.Item(137, 11).Value = 11249.0001
$ws.Cells.Item(137, 12).Value = 19243.8339
$ws.Cells.Item(137, 13).Value = -6149.000100000001
$ws.Cells.Item(137, 14).Value = -29443.8339

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4353925
$ws.Cells.Item(102, 9).Value = 5051326
$ws.Cells.Item(102, 11).Value = 5051326
$ws.Cells.Item(102, 13).Value = -5049704

$ws.Cells.Item(132, 8).Value = 6159.361
$ws.Cells.Item(132, 9).Value = 4624.1924
$ws.Cells.Item(132, 10).Value = 15029.223
$ws.Cells.Item(132, 11).Value = 13872.5772
$ws.Cells.Item(132, 12).Value = 45087.669
$ws.Cells.Item(132, 13).Value = -11342.5772
$ws.Cells.Item(132, 14).Value = -50147.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 2611
$ws.Cells.Item(3, 9).Value = 2222
$ws.Cells.Item(3, 11).Value = 2222
$ws.Cells.Item(3, 13).Value = -2110

$ws.Cells.Item(7, 8).Value = 5119.8037
$ws.Cells.Item(7, 9).Value = 4035.5806
$ws.Cells.Item(7, 11).Value = 4035.5806
$ws.Cells.Item(7, 13).Value = -3923.5806

$ws.Cells.Item(15, 8).Value = 2611
$ws.Cells.Item(15, 9).Value = 2222
$ws.Cells.Item(15, 11).Value = 2222
$ws.Cells.Item(15, 13).Value = -2052

$ws.Cells.Item(20, 8).Value = 7266.6665
$ws.Cells.Item(20, 9).Value = 6900
$ws.Cells.Item(20, 10).Value = 8000
$ws.Cells.Item(20, 11).Value = 6900
$ws.Cells.Item(20, 12).Value = 8000
$ws.Cells.Item(20, 13).Value = -6674
$ws.Cells.Item(20, 14).Value = -8452

$ws.Cells.Item(21, 8).Value = 1982.8334
$ws.Cells.Item(21, 10).Value = 1982.8334
$ws.Cells.Item(21, 12).Value = 1982.8334
$ws.Cells.Item(21, 14).Value = -2330.8334

$ws.Cells.Item(68, 8).Value = 1000
$ws.Cells.Item(68, 9).Value = 1000
$ws.Cells.Item(68, 11).Value = 1000
$ws.Cells.Item(68, 13).Value = -251

$ws.Cells.Item(71, 8).Value = 1000
$ws.Cells.Item(71, 9).Value = 1000
$ws.Cells.Item(71, 11).Value = 5000
$ws.Cells.Item(71, 13).Value = -1256

$ws.Cells.Item(126, 8).Value = 5119.8037
$ws.Cells.Item(126, 9).Value = 4035.5806
$ws.Cells.Item(126, 11).Value = 12106.7418
$ws.Cells.Item(126, 13).Value = -9636.7418

$ws.Cells.Item(132, 8).Value = 30643.584
$ws.Cells.Item(132, 9).Value = 32884.816
$ws.Cells.Item(132, 11).Value = 98654.448
$ws.Cells.Item(132, 13).Value = -96124.448

$ws.Cells.Item(136, 8).Value = 57163.684
$ws.Cells.Item(136, 9).Value = 72473.10000000001
$ws.Cells.Item(136, 10).Value = 7833.3335
$ws.Cells.Item(136, 11).Value = 217419.3
$ws.Cells.Item(136, 12).Value = 23500.0005
$ws.Cells.Item(136, 13).Value = -214869.3
$ws.Cells.Item(136, 14).Value = -28600.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 1000
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -886
$ws.Cells.Item(3, 14).ClearContents()

$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()

$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 3429.1428
$ws.Cells.Item(126, 9).Value = 3819.8
$ws.Cells.Item(126, 10).Value = 2452.5
$ws.Cells.Item(126, 11).Value = 11459.4
$ws.Cells.Item(126, 12).Value = 7357.5
$ws.Cells.Item(126, 13).Value = -8989.400000000001
$ws.Cells.Item(126, 14).Value = -12297.5

$ws.Cells.Item(132, 8).Value = 14448331
$ws.Cells.Item(132, 9).Value = 15879156
$ws.Cells.Item(132, 10).Value = 1570909.1
$ws.Cells.Item(132, 11).Value = 47637468
$ws.Cells.Item(132, 12).Value = 4712727.300000001
$ws.Cells.Item(132, 13).Value = -47634938
$ws.Cells.Item(132, 14).Value = -4717787.300000001
